# "created own file for setting headers of the xlsx file"
#
# Rename the sheet to "Forsendelser" and replace the sample data:
#   - "Fra"/"Til" (cols A/B) become free-text shipment codes instead of
#     plain numbers
#   - "Sendingsdato" (col C) becomes a real date value (dd/mm/yyyy) instead
#     of a shared-string label
#   - numeric columns (D weight, E/F/G dimensions) keep their values/format

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- sheet name ---
$ws.Name = "Forsendelser"

# --- column C: date format for header + all data cells ---
$ws.Columns("C").NumberFormat = "dd/mm/yyyy;@"

# --- row 2 ---
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "0024"
$ws.Range("A2").Style = "Normal"

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2500"
$ws.Range("B2").Style = "Normal"

$ws.Range("C2").Value = (Get-Date -Year 2020 -Month 12 -Day 18).Date

# --- row 3 ---
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "0024"
$ws.Range("A3").Style = "Normal"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "7075"
$ws.Range("B3").Style = "Normal"

$ws.Range("C3").Value = (Get-Date -Year 2020 -Month 12 -Day 19).Date

# --- row 4 ---
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "0024"
$ws.Range("A4").Style = "Normal"

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "4200"
$ws.Range("B4").Style = "Normal"

$ws.Range("C4").Value = (Get-Date -Year 2020 -Month 12 -Day 20).Date

# --- selection / cursor position ---
$ws.Range("C4").Select()

# --- print / page setup ---
$ws.PageSetup.Orientation = 1
